$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.686.35'
$ws.Range("E2").Value = '  -1.30%  '
$ws.Range("D3").Value = '2.908.58'
$ws.Range("E3").Value = '  -1.63%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '586.48'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.58%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.71'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.61%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.507'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.09%  '
$ws.Range("D9").Value = '2.907.99'
$ws.Range("E9").Value = '  -1.68%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.92'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.71%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.149'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.90%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.436'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.49%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000236'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.73%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.91'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.42%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.125'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.46%  '
$ws.Range("D16").Value = '3.387.94'
$ws.Range("E16").Value = '  -2.03%  '
$ws.Range("D17").Value = '61.718.22'
$ws.Range("E17").Value = '  -1.30%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.61'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.44%  '
$ws.Range("D19").Value = '2.908.83'
$ws.Range("E19").Value = '  -0.81%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '438.06'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.54%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.43'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.35%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.660'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.53%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.94'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.53%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '81.02'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.73%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.95'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.55%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.16'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -7.78%  '
$ws.Range("B27").Value = 'Fetch.AI'
$ws.Range("C27").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.07'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.35%  '
$ws.Range("B28").Value = 'Dai'
$ws.Range("C28").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000105'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +19.60%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.20'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.79%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.56'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.97%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.11'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.67%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.109'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.88%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.29%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '25.89'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.71%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.975'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.54%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.51'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.43%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.03'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.13%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '49.08'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.30%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.00'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.50%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.36'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.73%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.116'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.34%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.272'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.24%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '38.86'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.45%  '
$ws.Range("D45").Value = '2.691.11'
$ws.Range("E45").Value = '  -0.79%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '133.74'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.69%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0337'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.04%  '
$ws.Range("B48").Value = 'Bittensor'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '343.58'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.98%  '
$ws.Range("B49").Value = 'USDe'
$ws.Range("C49").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.00'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.103'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.55%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '22.32'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.12%  '
